$wb = $excel.ActiveWorkbook

# KN01
$ws = $wb.Worksheets.Item("KN01")
$ws.Range("I9").Value = '($3925) fr 2024 + $2641 Jan= $6566. , Dee pd $200'

# KN02
$ws = $wb.Worksheets.Item("KN02")
$ws.Range("B20").Value = 2775
$ws.Range("C20").Value = 'Fort Worth'
$ws.Range("D20").Value = 2775
$ws.Range("F20").Value = 2775
$ws.Range("I20").Value = ""

# KN03
$ws = $wb.Worksheets.Item("KN03")
$ws.Range("B4").Value = '2701 Echo Point Dr.'
$ws.Range("H4").Value = 'Kristine Zepeda'
$ws.Range("B9").Value = 5828
$ws.Range("C9").Value = 'Tarrant County'
$ws.Range("D9").Value = 1971
$ws.Range("E9").Value = 150
$ws.Range("F9").Value = 2121
$ws.Range("G9").Value = -3707
$ws.Range("H9").Value = -3707
$ws.Range("I9").Value = '$3128 (2024 under bal) + Rent $2700 = $5828 T pd $150'
$ws.Range("B10").Value = 2700
$ws.Range("C10").Value = 'Tarrant County'
$ws.Range("D10").Value = 1971
$ws.Range("F10").Value = 1971
$ws.Range("G10").Value = -729
$ws.Range("H10").Value = -4436
$ws.Range("I10").Value = ""
$ws.Range("B11").Value = 2700
$ws.Range("C11").Value = 'Tarrant County'
$ws.Range("D11").Value = 1971
$ws.Range("E11").Value = 500
$ws.Range("F11").Value = 2471
$ws.Range("G11").Value = -229
$ws.Range("H11").Value = -4665
$ws.Range("I11").Value = ""
$ws.Range("B12").Value = 2700
$ws.Range("C12").Value = 'Tarrant County'
$ws.Range("D12").Value = 1971
$ws.Range("E12").Value = 730
$ws.Range("F12").Value = 2701
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = -4664
$ws.Range("I12").Value = ""
$ws.Range("B13").Value = 2700
$ws.Range("C13").Value = 'Tarrant County'
$ws.Range("D13").Value = 1971
$ws.Range("F13").Value = 1971
$ws.Range("G13").Value = -729
$ws.Range("H13").Value = -5393
$ws.Range("I13").Value = ""
$ws.Range("B14").Value = 2700
$ws.Range("C14").Value = 'Tarrant County'
$ws.Range("D14").Value = 1971
$ws.Range("E14").Value = 200
$ws.Range("F14").Value = 2171
$ws.Range("G14").Value = -529
$ws.Range("H14").Value = -5922
$ws.Range("I14").Value = ""
$ws.Range("B15").Value = 2700
$ws.Range("C15").Value = 'Tarrant County'
$ws.Range("D15").Value = 1971
$ws.Range("E15").Value = 201.5
$ws.Range("F15").Value = 2172.5
$ws.Range("G15").Value = -527.5
$ws.Range("H15").Value = -6449.5
$ws.Range("I15").Value = ""
$ws.Range("B16").Value = 2700
$ws.Range("C16").Value = 'Tarrant County'
$ws.Range("D16").Value = 2510
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 2710
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -6439.5
$ws.Range("I16").Value = 'Eff 8/1 (HS: $2510 + T: $190) T pd $200 9/2'
$ws.Range("B17").Value = 2700
$ws.Range("C17").Value = 'Tarrant County'
$ws.Range("E17").Value = 190
$ws.Range("F17").Value = 190
$ws.Range("G17").Value = -2510
$ws.Range("H17").Value = -8949.5
$ws.Range("I17").Value = 'THS pd ($2510) 2 = $5020'
$ws.Range("B18").Value = 2700
$ws.Range("C18").Value = 'Tarrant County'
$ws.Range("D18").Value = 5020
$ws.Range("F18").Value = 5020
$ws.Range("G18").Value = 2320
$ws.Range("H18").Value = -6629.5
$ws.Range("I18").Value = 'THS pd ($2510) 2 = $5020'
$ws.Range("B19").Value = 2250
$ws.Range("C19").Value = 'Tarrant County'
$ws.Range("G19").Value = -2250
$ws.Range("H19").Value = -8879.5
$ws.Range("I19").Value = 'Housing Terminated contract Oct 31, 2025   25 days x $90 = $2,250 (Nov)'
$ws.Range("C20").Value = 'Tarrant County'
$ws.Range("H20").Value = -8879.5
$ws.Range("I20").Value = ""

# KN04
$ws = $wb.Worksheets.Item("KN04")
$ws.Range("B4").Value = '2604 Poplar Springs'
$ws.Range("H4").Value = 'Sandra Flory'
$ws.Range("B9").Value = 2782
$ws.Range("C9").Value = 'Fort Worth'
$ws.Range("D9").Value = 1546
$ws.Range("E9").Value = 1236
$ws.Range("F9").Value = 2782
$ws.Range("I9").Value = ""
$ws.Range("B10").Value = 2782
$ws.Range("C10").Value = 'Fort Worth'
$ws.Range("D10").Value = 1546
$ws.Range("E10").Value = 1236
$ws.Range("F10").Value = 2782
$ws.Range("I10").Value = ""
$ws.Range("B11").Value = 2782
$ws.Range("C11").Value = 'Fort Worth'
$ws.Range("D11").Value = 1546
$ws.Range("E11").Value = 1236
$ws.Range("F11").Value = 2782
$ws.Range("I11").Value = ""
$ws.Range("B12").Value = 2782
$ws.Range("C12").Value = 'Fort Worth'
$ws.Range("D12").Value = 1546
$ws.Range("E12").Value = 1236
$ws.Range("F12").Value = 2782
$ws.Range("I12").Value = 'T pd: $1000 4/2 + $236 $/4/3 = $1,236'
$ws.Range("B13").Value = 2782
$ws.Range("C13").Value = 'Fort Worth'
$ws.Range("D13").Value = 1546
$ws.Range("E13").Value = 1236
$ws.Range("F13").Value = 2782
$ws.Range("I13").Value = ""
$ws.Range("B14").Value = 2782
$ws.Range("C14").Value = 'Fort Worth'
$ws.Range("D14").Value = 1546
$ws.Range("E14").Value = 1236
$ws.Range("F14").Value = 2782
$ws.Range("I14").Value = ""
$ws.Range("B15").Value = 2782
$ws.Range("C15").Value = 'Fort Worth'
$ws.Range("D15").Value = 1546
$ws.Range("E15").Value = 1236
$ws.Range("F15").Value = 2782
$ws.Range("I15").Value = ""
$ws.Range("B16").Value = 2782
$ws.Range("C16").Value = 'Fort Worth'
$ws.Range("D16").Value = 1546
$ws.Range("E16").Value = 1236
$ws.Range("F16").Value = 2782
$ws.Range("I16").Value = ""
$ws.Range("B17").Value = 2782
$ws.Range("C17").Value = 'Fort Worth'
$ws.Range("D17").Value = 2008
$ws.Range("E17").Value = 774
$ws.Range("F17").Value = 2782
$ws.Range("I17").Value = 'Eff 9/1 (HS: $2,008 + T: $774), T pd $774'
$ws.Range("B18").Value = 2782
$ws.Range("C18").Value = 'Fort Worth'
$ws.Range("D18").Value = 2008
$ws.Range("E18").Value = 774
$ws.Range("F18").Value = 2782
$ws.Range("I18").Value = ""
$ws.Range("B19").Value = 2782
$ws.Range("C19").Value = 'Fort Worth'
$ws.Range("D19").Value = 2008
$ws.Range("E19").Value = 774
$ws.Range("F19").Value = 2782
$ws.Range("I19").Value = ""
$ws.Range("B20").Value = 2782
$ws.Range("C20").Value = 'Fort Worth'
$ws.Range("D20").Value = 2008
$ws.Range("E20").Value = 774
$ws.Range("F20").Value = 2782
$ws.Range("I20").Value = ""

# KN05
$ws = $wb.Worksheets.Item("KN05")
$ws.Range("B4").Value = '2817 Country Creek'
$ws.Range("H4").Value = 'Alisha Friddle'
$ws.Range("B9").Value = 2632
$ws.Range("C9").Value = 'Fort Worth'
$ws.Range("D9").Value = 2632
$ws.Range("F9").Value = 2632
$ws.Range("I9").Value = ""
$ws.Range("B10").Value = 2632
$ws.Range("C10").Value = 'Fort Worth'
$ws.Range("D10").Value = 2632
$ws.Range("F10").Value = 2632
$ws.Range("I10").Value = ""
$ws.Range("B11").Value = 2632
$ws.Range("C11").Value = 'Fort Worth'
$ws.Range("D11").Value = 2632
$ws.Range("F11").Value = 2632
$ws.Range("I11").Value = ""
$ws.Range("B12").Value = 2632
$ws.Range("C12").Value = 'Fort Worth'
$ws.Range("D12").Value = 2632
$ws.Range("F12").Value = 2632
$ws.Range("I12").Value = ""
$ws.Range("B13").Value = 2632
$ws.Range("C13").Value = 'Fort Worth'
$ws.Range("D13").Value = 2632
$ws.Range("F13").Value = 2632
$ws.Range("I13").Value = ""
$ws.Range("B14").Value = 2632
$ws.Range("C14").Value = 'Fort Worth'
$ws.Range("D14").Value = 2632
$ws.Range("F14").Value = 2632
$ws.Range("I14").Value = ""
$ws.Range("B15").Value = 2632
$ws.Range("C15").Value = 'Fort Worth'
$ws.Range("D15").Value = 2632
$ws.Range("F15").Value = 2632
$ws.Range("I15").Value = ""
$ws.Range("B16").Value = 2632
$ws.Range("C16").Value = 'Fort Worth'
$ws.Range("D16").Value = 2632
$ws.Range("F16").Value = 2632
$ws.Range("I16").Value = ""
$ws.Range("B17").Value = 2632
$ws.Range("C17").Value = 'Fort Worth'
$ws.Range("D17").Value = 2218
$ws.Range("E17").Value = 485
$ws.Range("F17").Value = 2703
$ws.Range("G17").Value = 71
$ws.Range("H17").Value = 71
$ws.Range("I17").Value = 'Eff 11/1 (HS: $2,147+ T: $485) , T pd $300 8/28 + $185 9/4 = $485'
$ws.Range("B18").Value = 2632
$ws.Range("C18").Value = 'Fort Worth'
$ws.Range("D18").Value = 2218
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 2418
$ws.Range("G18").Value = -214
$ws.Range("H18").Value = -143
$ws.Range("I18").Value = 'HS pd $2147 + $71 = $2218, T pd $200 10/27'
$ws.Range("B19").Value = 2632
$ws.Range("C19").Value = 'Fort Worth'
$ws.Range("D19").Value = 2147
$ws.Range("E19").Value = 300
$ws.Range("F19").Value = 2447
$ws.Range("G19").Value = -185
$ws.Range("H19").Value = -328
$ws.Range("I19").Value = 'T pd $300 10/31'
$ws.Range("B20").Value = 2632
$ws.Range("C20").Value = 'Fort Worth'
$ws.Range("D20").Value = 2147
$ws.Range("E20").Value = 400
$ws.Range("F20").Value = 2547
$ws.Range("G20").Value = -85
$ws.Range("H20").Value = -413
$ws.Range("I20").Value = 'T pd $400 12/1'

# KN06
$ws = $wb.Worksheets.Item("KN06")
$ws.Range("B4").Value = '3872 Country Ln'
$ws.Range("H4").Value = 'Kenesha Jackson'
$ws.Range("B9").Value = 2450
$ws.Range("C9").Value = 'Fort Worth'
$ws.Range("D9").Value = 2450
$ws.Range("F9").Value = 2450
$ws.Range("I9").Value = ""
$ws.Range("B10").Value = 2450
$ws.Range("C10").Value = 'Fort Worth'
$ws.Range("D10").Value = 2450
$ws.Range("F10").Value = 2450
$ws.Range("I10").Value = ""
$ws.Range("B11").Value = 2450
$ws.Range("C11").Value = 'Fort Worth'
$ws.Range("D11").Value = 2450
$ws.Range("F11").Value = 2450
$ws.Range("I11").Value = ""
$ws.Range("B12").Value = 2450
$ws.Range("C12").Value = 'Fort Worth'
$ws.Range("D12").Value = 2450
$ws.Range("F12").Value = 2450
$ws.Range("I12").Value = ""
$ws.Range("B13").Value = 2450
$ws.Range("C13").Value = 'Fort Worth'
$ws.Range("D13").Value = 2450
$ws.Range("F13").Value = 2450
$ws.Range("I13").Value = ""
$ws.Range("B14").Value = 2450
$ws.Range("C14").Value = 'Fort Worth'
$ws.Range("D14").Value = 2450
$ws.Range("F14").Value = 2450
$ws.Range("I14").Value = ""
$ws.Range("B15").Value = 2450
$ws.Range("C15").Value = 'Fort Worth'
$ws.Range("D15").Value = 2450
$ws.Range("F15").Value = 2450
$ws.Range("I15").Value = ""
$ws.Range("B16").Value = 2450
$ws.Range("C16").Value = 'Fort Worth'
$ws.Range("D16").Value = 2450
$ws.Range("F16").Value = 2450
$ws.Range("I16").Value = ""
$ws.Range("B17").Value = 2450
$ws.Range("C17").Value = 'Fort Worth'
$ws.Range("D17").Value = 2450
$ws.Range("F17").Value = 2450
$ws.Range("I17").Value = ""
$ws.Range("B18").Value = 2450
$ws.Range("C18").Value = 'Fort Worth'
$ws.Range("D18").Value = 2450
$ws.Range("F18").Value = 2450
$ws.Range("I18").Value = ""
$ws.Range("B19").Value = 2450
$ws.Range("C19").Value = 'Fort Worth'
$ws.Range("D19").Value = 1742
$ws.Range("E19").Value = 708
$ws.Range("F19").Value = 2450
$ws.Range("I19").Value = 'Eff 11/1 (HS $1742 + T $708) = $2,450, T pd $708 10/24'
$ws.Range("B20").Value = 2450
$ws.Range("C20").Value = 'Fort Worth'
$ws.Range("D20").Value = 1742
$ws.Range("E20").Value = 708
$ws.Range("F20").Value = 2450
$ws.Range("I20").Value = ""

# KN07
$ws = $wb.Worksheets.Item("KN07")
$ws.Range("B4").Value = '3841 Country Ln'
$ws.Range("H4").Value = 'Laquita Justice'
$ws.Range("B9").Value = 3612
$ws.Range("C9").Value = 'Fort Worth'
$ws.Range("D9").Value = 1548
$ws.Range("E9").Value = 1000
$ws.Range("F9").Value = 2548
$ws.Range("G9").Value = -1064
$ws.Range("H9").Value = -1064
$ws.Range("I9").Value = 'Jan $2,315 + $1297 (2024 bal due) = $3,612      A pd $600 + L pd $400 = $1000'
$ws.Range("B10").Value = 2315
$ws.Range("C10").Value = 'Fort Worth'
$ws.Range("D10").Value = 1548
$ws.Range("E10").Value = 700
$ws.Range("F10").Value = 2248
$ws.Range("G10").Value = -67
$ws.Range("H10").Value = -1131
$ws.Range("I10").Value = '$200 (A pd) + $500 (L pd) = $700'
$ws.Range("B11").Value = 2315
$ws.Range("C11").Value = 'Fort Worth'
$ws.Range("D11").Value = 1548
$ws.Range("E11").Value = 1500
$ws.Range("F11").Value = 3048
$ws.Range("G11").Value = 733
$ws.Range("H11").Value = -398
$ws.Range("I11").Value = '$850 (A pd) + $650 (L Pd) = $1,500'
$ws.Range("B12").Value = 2315
$ws.Range("C12").Value = 'Fort Worth'
$ws.Range("D12").Value = 887
$ws.Range("E12").Value = 1000
$ws.Range("F12").Value = 1887
$ws.Range("G12").Value = -428
$ws.Range("H12").Value = -826
$ws.Range("I12").Value = 'Eff 4/1 (HS: 887 + T: $1,428)= $2,315.            A: pd $700 + L: pd $300 = $1000'
$ws.Range("B13").Value = 2315
$ws.Range("C13").Value = 'Fort Worth'
$ws.Range("D13").Value = 887
$ws.Range("E13").Value = 700
$ws.Range("F13").Value = 1587
$ws.Range("G13").Value = -728
$ws.Range("H13").Value = -1554
$ws.Range("I13").Value = '$700 (A pd)'
$ws.Range("B14").Value = 2315
$ws.Range("C14").Value = 'Fort Worth'
$ws.Range("D14").Value = 887
$ws.Range("E14").Value = 600
$ws.Range("F14").Value = 1487
$ws.Range("G14").Value = -828
$ws.Range("H14").Value = -2382
$ws.Range("I14").Value = '$600 (A pd 3 x $200)'
$ws.Range("B15").Value = 2315
$ws.Range("C15").Value = 'Fort Worth'
$ws.Range("D15").Value = 887
$ws.Range("E15").Value = 480
$ws.Range("F15").Value = 1367
$ws.Range("G15").Value = -948
$ws.Range("H15").Value = -3330
$ws.Range("I15").Value = '$300 (A pd) + $180 (L pd) = $480'
$ws.Range("B16").Value = 2315
$ws.Range("C16").Value = 'Fort Worth'
$ws.Range("D16").Value = 887
$ws.Range("E16").Value = 1200
$ws.Range("F16").Value = 2087
$ws.Range("G16").Value = -228
$ws.Range("H16").Value = -3558
$ws.Range("I16").Value = '$200+$200 (A pd) + $300+$500 (L pd) = $1200'
$ws.Range("B17").Value = 2315
$ws.Range("C17").Value = 'Fort Worth'
$ws.Range("D17").Value = 887
$ws.Range("E17").Value = 1550
$ws.Range("F17").Value = 2437
$ws.Range("G17").Value = 122
$ws.Range("H17").Value = -3436
$ws.Range("I17").Value = 'A pd $300 9/2 + $500 9/16 + $250 9/29 + L pd $500 9/12 = Total pd $1,550'
$ws.Range("B18").Value = 2315
$ws.Range("C18").Value = 'Fort Worth'
$ws.Range("D18").Value = 887
$ws.Range("E18").Value = 1100
$ws.Range("F18").Value = 1987
$ws.Range("G18").Value = -328
$ws.Range("H18").Value = -3764
$ws.Range("I18").Value = 'A pd $100 9/30 + $300 10/20 +$300 10/27 = $700 + L pd $400 = $1,100'
$ws.Range("B19").Value = 2315
$ws.Range("C19").Value = 'Fort Worth'
$ws.Range("D19").Value = 1968
$ws.Range("E19").Value = 447
$ws.Range("F19").Value = 2415
$ws.Range("G19").Value = 100
$ws.Range("H19").Value = -3664
$ws.Range("I19").Value = 'Eff 11/1 (HS: $1968 + T: $347)= $2,315.                       A pd $347 11/4 + $100 11/19 = $447'
$ws.Range("B20").Value = 2315
$ws.Range("C20").Value = 'Fort Worth'
$ws.Range("D20").Value = 1968
$ws.Range("E20").Value = 340
$ws.Range("F20").Value = 2308
$ws.Range("G20").Value = -7
$ws.Range("H20").Value = -3671
$ws.Range("I20").Value = 'A pd $240 12/1 + $100 12/15 = $340 (2025 Bal due $3,671)'

# KN08
$ws = $wb.Worksheets.Item("KN08")
$ws.Range("B4").Value = '3908 Irish Setter Dr.'
$ws.Range("H4").Value = 'Gabrielle Owens'
$ws.Range("B9").Value = 2485
$ws.Range("C9").Value = 'Fort Worth'
$ws.Range("D9").Value = 2142
$ws.Range("E9").Value = 343
$ws.Range("F9").Value = 2485
$ws.Range("I9").Value = ""
$ws.Range("B10").Value = 2485
$ws.Range("C10").Value = 'Fort Worth'
$ws.Range("D10").Value = 2142
$ws.Range("E10").Value = 343
$ws.Range("F10").Value = 2485
$ws.Range("I10").Value = ""
$ws.Range("B11").Value = 2485
$ws.Range("C11").Value = 'Fort Worth'
$ws.Range("D11").Value = 2142
$ws.Range("E11").Value = 343
$ws.Range("F11").Value = 2485
$ws.Range("I11").Value = ""
$ws.Range("B12").Value = 2485
$ws.Range("C12").Value = 'Fort Worth'
$ws.Range("D12").Value = 2142
$ws.Range("E12").Value = 343
$ws.Range("F12").Value = 2485
$ws.Range("I12").Value = ""
$ws.Range("B13").Value = 2485
$ws.Range("C13").Value = 'Fort Worth'
$ws.Range("D13").Value = 2142
$ws.Range("E13").Value = 343
$ws.Range("F13").Value = 2485
$ws.Range("I13").Value = ""
$ws.Range("B14").Value = 2485
$ws.Range("C14").Value = 'Fort Worth'
$ws.Range("D14").Value = 2485
$ws.Range("F14").Value = 2485
$ws.Range("I14").Value = 'Eff 6/1 FWHS: $2485'
$ws.Range("B15").Value = 2485
$ws.Range("C15").Value = 'Fort Worth'
$ws.Range("D15").Value = 2485
$ws.Range("F15").Value = 2485
$ws.Range("I15").Value = ""
$ws.Range("B16").Value = 2485
$ws.Range("C16").Value = 'Fort Worth'
$ws.Range("D16").Value = 2485
$ws.Range("F16").Value = 2485
$ws.Range("I16").Value = ""
$ws.Range("B17").Value = 2485
$ws.Range("C17").Value = 'Fort Worth'
$ws.Range("D17").Value = 2485
$ws.Range("F17").Value = 2485
$ws.Range("I17").Value = ""
$ws.Range("B18").Value = 2485
$ws.Range("C18").Value = 'Fort Worth'
$ws.Range("D18").Value = 2485
$ws.Range("F18").Value = 2485
$ws.Range("I18").Value = ""
$ws.Range("B19").Value = 2485
$ws.Range("C19").Value = 'Fort Worth'
$ws.Range("D19").Value = 2485
$ws.Range("F19").Value = 2485
$ws.Range("I19").Value = ""
$ws.Range("B20").Value = 2485
$ws.Range("C20").Value = 'Fort Worth'
$ws.Range("D20").Value = 2485
$ws.Range("F20").Value = 2485
$ws.Range("I20").Value = ""

# KN09
$ws = $wb.Worksheets.Item("KN09")
$ws.Range("B9").Value = 2350
$ws.Range("C9").Value = 'Tarrant County'
$ws.Range("D9").Value = 2350
$ws.Range("F9").Value = 2350
$ws.Range("I9").Value = ""
$ws.Range("B10").Value = 2350
$ws.Range("C10").Value = 'Tarrant County'
$ws.Range("D10").Value = 2350
$ws.Range("F10").Value = 2350
$ws.Range("I10").Value = ""
$ws.Range("B11").Value = 2350
$ws.Range("C11").Value = 'Tarrant County'
$ws.Range("D11").Value = 2350
$ws.Range("F11").Value = 2350
$ws.Range("I11").Value = ""
$ws.Range("B12").Value = 2350
$ws.Range("C12").Value = 'Tarrant County'
$ws.Range("D12").Value = 2350
$ws.Range("F12").Value = 2350
$ws.Range("I12").Value = ""
